$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header labels / values (case fix on "UserName", typo fix on password)
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("B2").Value = "TestingPass123"

# Make the header row bold
$ws.Range("A1:B1").Font.Bold = $true

# Page setup: portrait A4
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
